$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 402.5
$ws.Range("I2").Value = 305.2
$ws.Range("K2").Value = 305.2
$ws.Range("M2").Value = -192.2

$ws.Range("H5").Value = 72.5
$ws.Range("I5").Value = 72.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 72.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 42.5
$ws.Range("N5").ClearContents()

$ws.Range("H43").Value = 1999.75
$ws.Range("I43").Value = 1999.6666
$ws.Range("K43").Value = 1999.6666
$ws.Range("M43").Value = -1930.6666

$ws.Range("H132").Value = 5097.2617
$ws.Range("I132").Value = 4849.472
$ws.Range("J132").Value = 6584
$ws.Range("K132").Value = 14548.416
$ws.Range("L132").Value = 19752
$ws.Range("M132").Value = -12018.416
$ws.Range("N132").Value = -24812

$ws.Range("H135").Value = 1631.3334
$ws.Range("I135").Value = 947
$ws.Range("K135").Value = 8523
$ws.Range("M135").Value = -5988

$ws.Range("H137").Value = 5753
$ws.Range("I137").Value = 4698
$ws.Range("J137").Value = 7042.4443
$ws.Range("K137").Value = 14094
$ws.Range("L137").Value = 21127.3329
$ws.Range("M137").Value = -11544
$ws.Range("N137").Value = -26227.3329

$ws.Range("H138").Value = 7165.625
$ws.Range("J138").Value = 6990.8076
$ws.Range("L138").Value = 20972.4228
$ws.Range("N138").Value = -31252.4228

$ws.Range("H141").Value = 4886.4595
$ws.Range("I141").Value = 4751.4287
$ws.Range("K141").Value = 14254.2861
$ws.Range("M141").Value = -9074.286100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16134.292
$ws.Range("I32").Value = 11687.217
$ws.Range("J32").Value = 69499.2
$ws.Range("K32").Value = 11687.217
$ws.Range("L32").Value = 69499.2
$ws.Range("M32").Value = -11400.217
$ws.Range("N32").Value = -70073.2

$ws.Range("H34").Value = 22999
$ws.Range("I34").Value = 22999
$ws.Range("K34").Value = 22999
$ws.Range("M34").Value = -22728

$ws.Range("H97").Value = 965.9666999999999
$ws.Range("I97").Value = 959.4545000000001
$ws.Range("J97").Value = 983.875
$ws.Range("K97").Value = 959.4545000000001
$ws.Range("L97").Value = 983.875
$ws.Range("M97").Value = -463.4545000000001
$ws.Range("N97").Value = -1975.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3282.8333
$ws.Range("I94").Value = 3648.125
$ws.Range("K94").Value = 3648.125
$ws.Range("M94").Value = -3197.125

$ws.Range("H134").Value = 4878.2856
$ws.Range("I134").Value = 4560.1
$ws.Range("K134").Value = 13680.3
$ws.Range("M134").Value = -11145.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 623218.5600000001
$ws.Range("I99").Value = 800138.1
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 800138.1
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -798640.1
$ws.Range("N99").Value = -6996

$ws.Range("H126").Value = 623218.5600000001
$ws.Range("I126").Value = 800138.1
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 2400414.3
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -2397944.3
$ws.Range("N126").Value = -16940

$ws.Range("H134").Value = 2817.111
$ws.Range("I134").Value = 2514.6326
$ws.Range("J134").Value = 5781.4
$ws.Range("K134").Value = 7543.8978
$ws.Range("L134").Value = 17344.2
$ws.Range("M134").Value = -5008.8978
$ws.Range("N134").Value = -22414.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 8153.222
$ws.Range("I141").Value = 8153.222
$ws.Range("K141").Value = 24459.666
$ws.Range("M141").Value = -19279.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 63330
$ws.Range("J63").Value = 63330
$ws.Range("L63").Value = 63330
$ws.Range("N63").Value = -64702

$ws.Range("H66").Value = 63330
$ws.Range("J66").Value = 63330
$ws.Range("L66").Value = 189990
$ws.Range("N66").Value = -196854

$ws.Range("H97").Value = 749.8570999999999
$ws.Range("I97").Value = 810.8823
$ws.Range("J97").Value = 490.5
$ws.Range("K97").Value = 810.8823
$ws.Range("L97").Value = 490.5
$ws.Range("M97").Value = -314.8823
$ws.Range("N97").Value = -1482.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 8600
$ws.Range("I5").Value = 3000
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -2887
$ws.Range("N5").Value = -10226

$ws.Range("H16").Value = 1762.8572
$ws.Range("I16").Value = 1297
$ws.Range("K16").Value = 1297
$ws.Range("M16").Value = -1127

$ws.Range("H20").Value = 1878999.9
$ws.Range("J20").Value = 1878999.9
$ws.Range("L20").Value = 1878999.9
$ws.Range("N20").Value = -1879451.9

$ws.Range("H22").Value = 1089.7
$ws.Range("I22").Value = 398.33334
$ws.Range("J22").Value = 1386
$ws.Range("K22").Value = 398.33334
$ws.Range("L22").Value = 1386
$ws.Range("M22").Value = -103.33334
$ws.Range("N22").Value = -1976

$ws.Range("H27").Value = 1089.7
$ws.Range("I27").Value = 398.33334
$ws.Range("J27").Value = 1386
$ws.Range("K27").Value = 398.33334
$ws.Range("L27").Value = 1386
$ws.Range("M27").Value = -291.33334
$ws.Range("N27").Value = -1600

$ws.Range("H68").Value = 7093.467
$ws.Range("I68").Value = 7084.769
$ws.Range("J68").Value = 7150
$ws.Range("K68").Value = 7084.769
$ws.Range("L68").Value = 7150
$ws.Range("M68").Value = -6335.769
$ws.Range("N68").Value = -8648

$ws.Range("H71").Value = 7093.467
$ws.Range("I71").Value = 7084.769
$ws.Range("J71").Value = 7150
$ws.Range("K71").Value = 35423.845
$ws.Range("L71").Value = 7150
$ws.Range("M71").Value = -31679.845
$ws.Range("N71").Value = -43238

$ws.Range("H100").Value = 1979.2
$ws.Range("I100").Value = 1979.2
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1979.2
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1438.2
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H30").Value = 715
$ws.Range("I30").Value = 715
$ws.Range("K30").Value = 715
$ws.Range("M30").Value = -608

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H46").Value = 61884
$ws.Range("J46").Value = 61884
$ws.Range("L46").Value = 61884
$ws.Range("N46").Value = -62346

$ws.Range("H62").Value = 4386653
$ws.Range("J62").Value = 8300.6
$ws.Range("L62").Value = 8300.6
$ws.Range("N62").Value = -9548.6

$ws.Range("H65").Value = 4386653
$ws.Range("J65").Value = 8300.6
$ws.Range("L65").Value = 41503
$ws.Range("N65").Value = -47743

$ws.Range("H133").Value = 36189.668
$ws.Range("J133").Value = 35713.375
$ws.Range("L133").Value = 35713.375
$ws.Range("N133").Value = -45833.375

$ws.Range("H134").Value = 61884
$ws.Range("J134").Value = 61884
$ws.Range("L134").Value = 185652
$ws.Range("N134").Value = -190722

$ws.Range("H136").Value = 11529.913
$ws.Range("I136").Value = 11532.857
$ws.Range("K136").Value = 34598.571
$ws.Range("M136").Value = -32048.571
